$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 11:35"

# Row 18 - Belgica: updated case counts
$ws.Range("B18").Value = 54288
$ws.Range("C18").Value = 307
$ws.Range("D18").Value = 14111
$ws.Range("E18").Value = 31274
$ws.Range("F18").Value = 407
$ws.Range("G18").Value = 60
$ws.Range("H18").Value = 8903

# Row 33 - Banglades: updated case counts
$ws.Range("B33").Value = 18863
$ws.Range("C33").Value = 1041
$ws.Range("E33").Value = 15219
$ws.Range("G33").Value = 14
$ws.Range("H33").Value = 283

# Row 36 - Israel: updated case counts
$ws.Range("B36").Value = 16567
$ws.Range("C36").Value = 19
$ws.Range("D36").Value = 12364
$ws.Range("E36").Value = 3939
$ws.Range("F36").Value = 62

# New country "Indonesia" is inserted into the ranking right after "Japon" (row 37),
# pushing Rumania and Austria down one row each. Re-write rows 38-40 with the
# new country order and their updated case counts.

# Row 38 becomes Indonesia (was Rumania)
$ws.Range("A38").Value = "Indonesia"
$ws.Range("B38").Value = 16006
$ws.Range("C38").Value = 568
$ws.Range("D38").Value = 3518
$ws.Range("E38").Value = 11445
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 15
$ws.Range("H38").Value = 1043

# Row 39 becomes Rumania (was Austria)
$ws.Range("A39").Value = "Rumania"
$ws.Range("B39").Value = 16002
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 7961
$ws.Range("E39").Value = 6995
$ws.Range("F39").Value = 228
$ws.Range("G39").Value = 10
$ws.Range("H39").Value = 1046

# Row 40 becomes Austria (was Indonesia)
$ws.Range("A40").Value = "Austria"
$ws.Range("B40").Value = 15997
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 14304
$ws.Range("E40").Value = 1069
$ws.Range("F40").Value = 55
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 624

# Row 142 - Nepal: updated case counts
$ws.Range("B142").Value = 246
$ws.Range("C142").Value = 3
$ws.Range("E142").Value = 211
